# The SSO server has been deprecated; replace the "SSO" worksheet section
# (SSO / server / app_id / app_secret / google client id / client secret)
# with a new, shorter "Authentication" section.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the four now-obsolete detail rows (app_id, app_secret, google client
# id, client secret). Deleting whole rows shifts every following row up by
# four, which is what re-aligns the rest of the sheet with the new layout.
$ws.Rows("17:20").Delete()

# A14/A15/A16 currently hold the old "SSO" / "server" rows (row A17 is the
# blank separator row that used to precede "Blob storage"). Overwrite them
# with the new Authentication section's content.
$ws.Range("A14").Value = "Authentication"
$ws.Range("A14").Font.Bold = $true

$ws.Range("A15").Value = "Authentication provider"
$ws.Range("A15").Font.Bold = $false
$ws.Range("B15").Value = "Google / OpenIDConnect / PAM / LDAP"
$ws.Range("B15").Font.Bold = $false

$ws.Range("A16").Value = "credentials"
$ws.Range("A16").Font.Bold = $false
